$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 18:55"

# Row 4
$ws.Range("B4").Value = 7613291
$ws.Range("C4").Value = 11514
$ws.Range("D4").Value = 4831029
$ws.Range("E4").Value = 2567900
$ws.Range("G4").Value = 85
$ws.Range("H4").Value = 214362

# Row 22
$ws.Range("B22").Value = 324443
$ws.Range("C22").Value = 1429
$ws.Range("D22").Value = 285050
$ws.Range("E22").Value = 30952
$ws.Range("G22").Value = 57
$ws.Range("H22").Value = 8441

# Row 26
$ws.Range("B26").Value = 300850
$ws.Range("C26").Value = 822
$ws.Range("E26").Value = 29351
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 9599

# Row 29
$ws.Range("B29").Value = 166116
$ws.Range("C29").Value = 1645
$ws.Range("D29").Value = 140230
$ws.Range("E29").Value = 16406
$ws.Range("G29").Value = 18
$ws.Range("H29").Value = 9480

# Row 50
$ws.Range("B50").Value = 81711
$ws.Range("C50").Value = 1106
$ws.Range("D50").Value = 43991
$ws.Range("E50").Value = 36998
$ws.Range("G50").Value = 11
$ws.Range("H50").Value = 722

# Row 68
$ws.Range("B68").Value = 44482
$ws.Range("C68").Value = 988
$ws.Range("D68").Value = 19782
$ws.Range("E68").Value = 24294
$ws.Range("G68").Value = 8
$ws.Range("H68").Value = 406

# Row 71
$ws.Range("B71").Value = 40691
$ws.Range("C71").Value = 130
$ws.Range("D71").Value = 38496
$ws.Range("E71").Value = 1599
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 596

# Row 87
$ws.Range("B87").Value = 19842
$ws.Range("C87").Value = 229
$ws.Range("E87").Value = 9444
$ws.Range("G87").Value = 4
$ws.Range("H87").Value = 409

# Row 101
$ws.Range("B101").Value = 11626
$ws.Range("C101").Value = 54
$ws.Range("D101").Value = 9429
$ws.Range("E101").Value = 2074

# Row 108
$ws.Range("B108").Value = 9196
$ws.Range("C108").Value = 147
$ws.Range("D108").Value = 5852
$ws.Range("E108").Value = 3278
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 66

# Row 133
$ws.Range("B133").Value = 4741
$ws.Range("C133").Value = 26
$ws.Range("D133").Value = 2843
$ws.Range("E133").Value = 1817
$ws.Range("G133").Value = 2
$ws.Range("H133").Value = 81

# Row 145
$ws.Range("A145").Value = "Mali"
$ws.Range("B145").Value = 3184
$ws.Range("C145").Value = 14
$ws.Range("D145").Value = 2479
$ws.Range("E145").Value = 574
$ws.Range("H145").Value = 131

# Row 146
$ws.Range("A146").Value = "Botsuana"
$ws.Range("B146").Value = 3172
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 710
$ws.Range("E146").Value = 2446
$ws.Range("H146").Value = 16

# Row 195
$ws.Range("A195").Value = "Liechtenstein"
$ws.Range("B195").Value = 126
$ws.Range("C195").Value = 3
$ws.Range("D195").Value = 116
$ws.Range("E195").Value = 9

# Row 196
$ws.Range("A196").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B196").Value = 124
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 46
$ws.Range("E196").Value = 77

# Row 215
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# Row 216
$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

